$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 136 (id 134) ---
$ws.Range("A135").Copy($ws.Range("A136"))
$ws.Range("A136").Value = 134
$ws.Range("B136").Value = 7011630
$ws.Range("C136").Value = "Azerbaijan Premier League"
$ws.Range("D136").Value = "Azerbaijan Premier League"
$ws.Range("E135").Copy($ws.Range("E136"))
$ws.Range("E136").Value = 45381.39583333334
$ws.Range("F136").Value = "Sabail FC"
$ws.Range("G136").Value = "PFK Turan Tovuz"
$ws.Range("K136").Value = 2.3
$ws.Range("L136").Value = 3.2
$ws.Range("M136").Value = 2.75
$ws.Range("N136").Value = 2.3
$ws.Range("O136").Value = 3.2
$ws.Range("P136").Value = 2.75
$ws.Range("Q136").Value = -0.25
$ws.Range("R136").Value = 2.05
$ws.Range("S136").Value = 1.75
$ws.Range("T136").Value = 2.25
$ws.Range("U136").Value = 1.875
$ws.Range("V136").Value = 1.925
$ws.Range("W136").Value = 0
$ws.Range("X136").Value = 0
$ws.Range("Y136").Value = 0
$ws.Range("Z136").Value = 0
$ws.Range("AA136").Value = 0

# --- Row 137 (id 135) ---
$ws.Range("A135").Copy($ws.Range("A137"))
$ws.Range("A137").Value = 135
$ws.Range("B137").Value = 7011629
$ws.Range("C137").Value = "Azerbaijan Premier League"
$ws.Range("D137").Value = "Azerbaijan Premier League"
$ws.Range("E135").Copy($ws.Range("E137"))
$ws.Range("E137").Value = 45381.5
$ws.Range("F137").Value = "Sabah"
$ws.Range("G137").Value = "Zira IK"
$ws.Range("K137").Value = 2.1
$ws.Range("L137").Value = 3.2
$ws.Range("M137").Value = 3.1
$ws.Range("N137").Value = 2.15
$ws.Range("O137").Value = 3.2
$ws.Range("P137").Value = 3.1
$ws.Range("Q137").Value = -0.25
$ws.Range("R137").Value = 1.9
$ws.Range("S137").Value = 1.9
$ws.Range("T137").Value = 2.25
$ws.Range("U137").Value = 1.925
$ws.Range("V137").Value = 1.875
$ws.Range("W137").Value = 0
$ws.Range("X137").Value = 0
$ws.Range("Y137").Value = 0
$ws.Range("Z137").Value = 0
$ws.Range("AA137").Value = 0

Write-Output "applied"
